# "Forgot to order table" - sort the per-language rows (rows 2-20) into the
# intended display order while keeping each row's values attached to its language.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing table (language label + the 5 numeric columns) before overwriting it.
$rows = @{}
for ($r = 2; $r -le 20; $r++) {
    $lang = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    $d = $ws.Cells.Item($r, 4).Value()
    $e = $ws.Cells.Item($r, 5).Value()
    $f = $ws.Cells.Item($r, 6).Value()
    $rows[$lang] = @($b, $c, $d, $e, $f)
}

# Desired language order (alphabetised-by-language-family grouping used in the paper table).
$newOrder = @(
    "German", "Spanish", "Slovak", "Norwegian", "Greek", "Chinese", "Vietnamese", "Thai", "Cantonese", "Indonesian", "Finnish", "Basque", "Korean", "Japanese", "Turkish", "Arabic", "Hebrew", "Algerian", "Maltese"
)

# Write the rows back out in the new order, carrying each language's own values with it.
$r = 2
foreach ($lang in $newOrder) {
    $vals = $rows[$lang]
    $ws.Cells.Item($r, 1).Value = $lang
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $r++
}
